# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split runs that currently end with a trailing space into two runs:
# the bare word, and a separate run containing just the space. This
# restores the pre-consolidation run layout (one run per "word token",
# with standalone single-space runs in between) without changing any
# visible text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 ("Title 1"): "Testing custom properties" ---------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
# "Testing custom properties"
#  123456789012345678901234567
#          ^8      ^15
$title.Characters(8, 1).Text = " "
$title.Characters(15, 1).Text = " "

# --- Shape 2 ("Subtitle 2"): "This is a subtitle<br><br>A. M." ------
$subtitle = $s.Shapes.Item(2).TextFrame.TextRange
# "This is a subtitle" then two line breaks (each counts as 1 char)
# then "A. M."
# T(1)h(2)i(3)s(4) (5)i(6)s(7) (8)a(9) (10)s(11)u(12)b(13)t(14)i(15)t(16)l(17)e(18)
# <br>(19) <br>(20) A(21).(22) (23)M(24).(25)
$subtitle.Characters(5, 1).Text = " "
$subtitle.Characters(8, 1).Text = " "
$subtitle.Characters(10, 1).Text = " "
$subtitle.Characters(23, 1).Text = " "
